{"js": "// Replace every occurrence of the old Slovak sentence with the new,\n// re-worded one: \"V roku S\u00fahvezdie Lev 2022: ...\" ->\n// \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Lev: ...\"\nconst oldText = \"V roku S\u00fahvezdie Lev 2022: 14. \u2013 23. apr\u00edla, 14. \u2013 23. m\u00e1ja\";\nconst newText = \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Lev: 14. \u2013 23. apr\u00edla, 14. \u2013 23. m\u00e1ja\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace every occurrence of the old Slovak sentence with the new,\n# re-worded one: \"V roku S\u00fahvezdie Lev 2022: ...\" ->\n# \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Lev: ...\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"V roku S\u00fahvezdie Lev 2022: 14. \u2013 23. apr\u00edla, 14. \u2013 23. m\u00e1ja\"\n$find.Replacement.Text = \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Lev: 14. \u2013 23. apr\u00edla, 14. \u2013 23. m\u00e1ja\"\n\n$find.Execute(\n    $find.Text,               # FindText\n    $false,                   # MatchCase\n    $false,                   # MatchWholeWord\n    $false,                   # MatchWildcards\n    $false,                   # MatchSoundsLike\n    $false,                   # MatchAllWordForms\n    $true,                    # Forward\n    1,                        # Wrap (wdFindContinue)\n    $false,                   # Format\n    $find.Replacement.Text,   # ReplaceWith\n    2                         # Replace (wdReplaceAll)\n)\n"}
